$wb = $excel.ActiveWorkbook

# This script re-applies a batch of refreshed market-price / profit values
# (columns H-N) across several Leve-profit sheets, as produced by the
# scheduled data-refresh runner. Values are written cell-by-cell via
# Range.Value so any existing formatting/number-format is preserved.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1863
$ws.Range("I58").Value = 1657.5
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 4972.5
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -4822.5
$ws.Range("N58").Value = -6300
$ws.Range("H62").Value = 12877.637
$ws.Range("I62").Value = 24766.223
$ws.Range("J62").Value = 4647.077
$ws.Range("K62").Value = 24766.223
$ws.Range("L62").Value = 4647.077
$ws.Range("M62").Value = -24142.223
$ws.Range("N62").Value = -5895.077
$ws.Range("H65").Value = 12877.637
$ws.Range("I65").Value = 24766.223
$ws.Range("J65").Value = 4647.077
$ws.Range("K65").Value = 123831.115
$ws.Range("L65").Value = 23235.385
$ws.Range("M65").Value = -120711.115
$ws.Range("N65").Value = -29475.385
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480
$ws.Range("H111").Value = 314.83334
$ws.Range("I111").Value = 296.33334
$ws.Range("J111").Value = 333.33334
$ws.Range("K111").Value = 889.0000200000001
$ws.Range("L111").Value = 1000.00002
$ws.Range("M111").Value = 2177.99998
$ws.Range("N111").Value = -7134.00002
$ws.Range("H135").Value = 1172.8182
$ws.Range("I135").Value = 529.5
$ws.Range("J135").Value = 1944.8
$ws.Range("K135").Value = 4765.5
$ws.Range("L135").Value = 17503.2
$ws.Range("M135").Value = -2230.5
$ws.Range("N135").Value = -22573.2
$ws.Range("H137").Value = 3369.6775
$ws.Range("I137").Value = 3402.4
$ws.Range("J137").Value = 3233.3333
$ws.Range("K137").Value = 10207.2
$ws.Range("L137").Value = 9699.999899999999
$ws.Range("M137").Value = -7657.200000000001
$ws.Range("N137").Value = -14799.9999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 272.375
$ws.Range("I5").Value = 297.14285
$ws.Range("J5").Value = 99
$ws.Range("K5").Value = 297.14285
$ws.Range("L5").Value = 99
$ws.Range("M5").Value = -185.14285
$ws.Range("N5").Value = -323
$ws.Range("H28").Value = 20932.5
$ws.Range("I28").Value = 11388.2
$ws.Range("J28").Value = 36839.668
$ws.Range("K28").Value = 11388.2
$ws.Range("L28").Value = 36839.668
$ws.Range("M28").Value = -11196.2
$ws.Range("N28").Value = -37223.668
$ws.Range("H99").Value = 20932.5
$ws.Range("I99").Value = 11388.2
$ws.Range("J99").Value = 36839.668
$ws.Range("K99").Value = 11388.2
$ws.Range("L99").Value = 36839.668
$ws.Range("M99").Value = -8393.200000000001
$ws.Range("N99").Value = -42829.668

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 272.375
$ws.Range("I4").Value = 297.14285
$ws.Range("J4").Value = 99
$ws.Range("K4").Value = 297.14285
$ws.Range("L4").Value = 99
$ws.Range("M4").Value = -182.14285
$ws.Range("N4").Value = -329
$ws.Range("H105").Value = 1807.1621
$ws.Range("I105").Value = 1720.2084
$ws.Range("J105").Value = 1967.6923
$ws.Range("K105").Value = 1720.2084
$ws.Range("L105").Value = 1967.6923
$ws.Range("M105").Value = 26.79160000000002
$ws.Range("N105").Value = -5461.6923

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 269.7143
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 502.66666
$ws.Range("K5").Value = 95
$ws.Range("L5").Value = 502.66666
$ws.Range("M5").Value = 17
$ws.Range("N5").Value = -726.66666
$ws.Range("H25").Value = 4002.75
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10348
$ws.Range("H134").Value = 2213.9355
$ws.Range("I134").Value = 1488.2667
$ws.Range("J134").Value = 2894.25
$ws.Range("K134").Value = 4464.800099999999
$ws.Range("L134").Value = 8682.75
$ws.Range("M134").Value = -1929.800099999999
$ws.Range("N134").Value = -13752.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1686.1609
$ws.Range("J131").Value = 1302.0494
$ws.Range("L131").Value = 3906.148200000001
$ws.Range("N131").Value = -13986.1482

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 16692970
$ws.Range("I18").Value = 25004452
$ws.Range("K18").Value = 25004452
$ws.Range("M18").Value = -25004159
$ws.Range("H21").Value = 4317
$ws.Range("I21").Value = 4475.5
$ws.Range("J21").Value = 4000
$ws.Range("K21").Value = 4475.5
$ws.Range("L21").Value = 4000
$ws.Range("M21").Value = -4302.5
$ws.Range("N21").Value = -4346
$ws.Range("H30").Value = 4317
$ws.Range("I30").Value = 4475.5
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 4475.5
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = -4370.5
$ws.Range("N30").Value = -4210
$ws.Range("H80").Value = 5483.3335
$ws.Range("I80").Value = 7761.5386
$ws.Range("K80").Value = 7761.5386
$ws.Range("M80").Value = -6763.5386
$ws.Range("H83").Value = 5483.3335
$ws.Range("I83").Value = 7761.5386
$ws.Range("K83").Value = 38807.693
$ws.Range("M83").Value = -33815.693
$ws.Range("H97").Value = 1017.53845
$ws.Range("I97").Value = 935.6667
$ws.Range("K97").Value = 935.6667
$ws.Range("M97").Value = -439.6667
$ws.Range("H132").Value = 4574
$ws.Range("I132").Value = 6901.091
$ws.Range("J132").Value = 3355.0476
$ws.Range("K132").Value = 20703.273
$ws.Range("L132").Value = 10065.1428
$ws.Range("M132").Value = -18173.273
$ws.Range("N132").Value = -15125.1428

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2251.6365
$ws.Range("I68").Value = 1818
$ws.Range("J68").Value = 2613
$ws.Range("K68").Value = 1818
$ws.Range("L68").Value = 2613
$ws.Range("M68").Value = -1069
$ws.Range("N68").Value = -4111
$ws.Range("H71").Value = 2251.6365
$ws.Range("I71").Value = 1818
$ws.Range("J71").Value = 2613
$ws.Range("K71").Value = 9090
$ws.Range("L71").Value = 13065
$ws.Range("M71").Value = -5346
$ws.Range("N71").Value = -20553
$ws.Range("H82").Value = 1985.9445
$ws.Range("J82").Value = 1896.6875
$ws.Range("L82").Value = 1896.6875
$ws.Range("N82").Value = -2618.6875
$ws.Range("H85").Value = 1985.9445
$ws.Range("J85").Value = 1896.6875
$ws.Range("L85").Value = 1896.6875
$ws.Range("N85").Value = -4392.6875
$ws.Range("H100").Value = 45458416
$ws.Range("I100").Value = 5270.923
$ws.Range("J100").Value = 111112950
$ws.Range("K100").Value = 5270.923
$ws.Range("L100").Value = 111112950
$ws.Range("M100").Value = -4729.923
$ws.Range("N100").Value = -111114032
$ws.Range("H136").Value = 4195.185
$ws.Range("I136").Value = 2617.0344
$ws.Range("J136").Value = 6025.84
$ws.Range("K136").Value = 7851.1032
$ws.Range("L136").Value = 18077.52
$ws.Range("M136").Value = -5301.1032
$ws.Range("N136").Value = -23177.52

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 4514
$ws.Range("I38").Value = 5056
$ws.Range("K38").Value = 5056
$ws.Range("M38").Value = -4583
$ws.Range("H48").Value = 9975
$ws.Range("J48").Value = 9975
$ws.Range("L48").Value = 9975
$ws.Range("N48").Value = -11113
$ws.Range("H49").Value = 727486.5600000001
$ws.Range("I49").Value = 1671818.6
$ws.Range("J49").Value = 19237.5
$ws.Range("K49").Value = 1671818.6
$ws.Range("L49").Value = 19237.5
$ws.Range("M49").Value = -1671588.6
$ws.Range("N49").Value = -19697.5
$ws.Range("H62").Value = 4530
$ws.Range("I62").Value = 4990
$ws.Range("J62").Value = 4464.2856
$ws.Range("K62").Value = 4990
$ws.Range("L62").Value = 4464.2856
$ws.Range("M62").Value = -4366
$ws.Range("N62").Value = -5712.2856
$ws.Range("H65").Value = 4530
$ws.Range("I65").Value = 4990
$ws.Range("J65").Value = 4464.2856
$ws.Range("K65").Value = 24950
$ws.Range("L65").Value = 22321.428
$ws.Range("M65").Value = -21830
$ws.Range("N65").Value = -28561.428
$ws.Range("H81").Value = 1929.1666
$ws.Range("I81").Value = 1485
$ws.Range("J81").Value = 4150
$ws.Range("K81").Value = 2970
$ws.Range("L81").Value = 8300
$ws.Range("M81").Value = -1909
$ws.Range("N81").Value = -10422
$ws.Range("H84").Value = 1929.1666
$ws.Range("I84").Value = 1485
$ws.Range("J84").Value = 4150
$ws.Range("K84").Value = 14850
$ws.Range("L84").Value = 41500
$ws.Range("M84").Value = -9546
$ws.Range("N84").Value = -52108
$ws.Range("H100").Value = 9307.375
$ws.Range("I100").Value = 5196
$ws.Range("J100").Value = 13418.75
$ws.Range("K100").Value = 10392
$ws.Range("L100").Value = 26837.5
$ws.Range("M100").Value = -9851
$ws.Range("N100").Value = -27919.5
$ws.Range("H136").Value = 23411364
$ws.Range("I136").Value = 29442636
$ws.Range("J136").Value = 14495571
$ws.Range("K136").Value = 88327908
$ws.Range("L136").Value = 43486713
$ws.Range("M136").Value = -88325358
$ws.Range("N136").Value = -43491813

